$d = $word.ActiveDocument

# Paragraph indices 74-83 (1-based, Word Paragraphs collection) get bold
# paragraph-mark + run formatting; paragraph 83 additionally gets its last
# run split around "imagen alternado" bracketed by proofErr gram tags.
$targets = 74..83

$xml0 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="06AFC690" w14:textId="77777777" w:rsidR="002D7A5B" w:rsidRPr="002D7A5B" w:rsidRDefault="002D7A5B" w:rsidP="002D7A5B"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">3. Completamente responsive </w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/><w:b/><w:bCs/></w:rPr><w:t>⚠️</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> Estructura responsive presente, pero algunos order-md-* están duplicados o mal puestos. Ajustar orden de columnas para cumplir los “3 cambios de orden”.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(74).Range.InsertXML($xml0)

$xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="20613DF5" w14:textId="77777777" w:rsidR="002D7A5B" w:rsidRPr="002D7A5B" w:rsidRDefault="002D7A5B" w:rsidP="002D7A5B"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">4. Estilos coherentes y resultado profesional </w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/><w:b/><w:bCs/></w:rPr><w:t>⚠️</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> Bien enfocado; faltaría pulir pequeños detalles (alineación, márgenes, consistencia de colores).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(75).Range.InsertXML($xml1)

$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4D92A9F5" w14:textId="77777777" w:rsidR="002D7A5B" w:rsidRPr="002D7A5B" w:rsidRDefault="002D7A5B" w:rsidP="002D7A5B"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">9. Botones en la página </w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/><w:b/><w:bCs/></w:rPr><w:t>❌</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> No hay ning</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/></w:rPr><w:t>ú</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>n bot</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/></w:rPr><w:t>ó</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>n visible ahora (los eliminaste en la secci</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/></w:rPr><w:t>ó</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>n de informaci</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/></w:rPr><w:t>ó</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>n pr</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/></w:rPr><w:t>á</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>ctica). A</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/></w:rPr><w:t>ñ</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">adir al menos 2 botones (ej. </w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/></w:rPr><w:t>“</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>C</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/></w:rPr><w:t>ó</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>mo llegar</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/></w:rPr><w:t>”</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/></w:rPr><w:t>“</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Ver m</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/></w:rPr><w:t>á</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>s</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/></w:rPr><w:t>”</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/></w:rPr><w:t>“</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Reservar alojamiento</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/></w:rPr><w:t>”</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(76).Range.InsertXML($xml2)

$xml3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0E8F0090" w14:textId="77777777" w:rsidR="002D7A5B" w:rsidRPr="002D7A5B" w:rsidRDefault="002D7A5B" w:rsidP="002D7A5B"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">10. Tabla responsive con bordes y striped </w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/><w:b/><w:bCs/></w:rPr><w:t>❌</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> Falta totalmente la tabla. Debes a</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/></w:rPr><w:t>ñ</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>adir &lt;div class="table-responsive"&gt;&lt;table class="table table-striped table-bordered"&gt;</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/></w:rPr><w:t>…</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>&lt;/table&gt;&lt;/div&gt; en alguna secci</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/></w:rPr><w:t>ó</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">n, por ejemplo </w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/></w:rPr><w:t>“</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Informaci</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/></w:rPr><w:t>ó</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>n pr</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/></w:rPr><w:t>á</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>ctica</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/></w:rPr><w:t>”</w:t></w:r><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(77).Range.InsertXML($xml3)

$xml4 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="1F153EA6" w14:textId="77777777" w:rsidR="002D7A5B" w:rsidRPr="002D7A5B" w:rsidRDefault="002D7A5B" w:rsidP="002D7A5B"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>11. Carrusel de Bootstrap con 3 imágenes y texto</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(78).Range.InsertXML($xml4)

$xml5 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0D1131D6" w14:textId="77777777" w:rsidR="002D7A5B" w:rsidRPr="002D7A5B" w:rsidRDefault="002D7A5B" w:rsidP="002D7A5B"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/><w:b/><w:bCs/></w:rPr><w:t>⚠️</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(79).Range.InsertXML($xml5)

$xml6 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="67B195D0" w14:textId="77777777" w:rsidR="002D7A5B" w:rsidRPr="002D7A5B" w:rsidRDefault="002D7A5B" w:rsidP="002D7A5B"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Carrusel bien hecho, pero tienes las 3 diapositivas con class="carousel-item active". Solo la primera debe ser active. Corrige eso para que funcione correctamente.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(80).Range.InsertXML($xml6)

$xml7 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="22CD39CF" w14:textId="77777777" w:rsidR="002D7A5B" w:rsidRPr="002D7A5B" w:rsidRDefault="002D7A5B" w:rsidP="002D7A5B"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>14. Columnas que cambian de disposición según versión</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(81).Range.InsertXML($xml7)

$xml8 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="362F8D98" w14:textId="77777777" w:rsidR="002D7A5B" w:rsidRPr="002D7A5B" w:rsidRDefault="002D7A5B" w:rsidP="002D7A5B"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="002D7A5B"><w:rPr><w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/><w:b/><w:bCs/></w:rPr><w:t>⚠️</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(82).Range.InsertXML($xml8)

$xml9 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="5957E06D" w14:textId="77777777" w:rsidR="002D7A5B" w:rsidRPr="002D7A5B" w:rsidRDefault="002D7A5B" w:rsidP="002D7A5B"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Las columnas usan order-md-*, pero hay errores (duplicados o en orden inverso). Corrige al menos 3 cambios reales de orden (ej. texto–</w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>imagen alternado</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="002D7A5B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(83).Range.InsertXML($xml9)
